$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data records currently stored in rows 13 and 14 were swapped
# (row 13 now holds the record that used to be in row 14, and vice
# versa). Re-create that by swapping every column's content between
# the two rows.

# Columns that genuinely hold numbers / booleans in this sheet - every
# other used column is free-form text (some of which, like "Antal" /
# column I, happens to look numeric, e.g. "1", "5").
$numericCols = @(1, 2, 5, 17, 18, 19)        # A, B, E, Q, R, S
$booleanCols = @(30, 31, 33)                 # AD, AE, AG

$firstCol = 1           # A
$lastCol  = 51          # AY

function Set-SwappedValue($targetCell, $newValue, $isText) {
    $isBlank = ($newValue -eq $null) -or ($newValue -eq "")

    if ($isText -and -not $isBlank) {
        # Force text interpretation first so digit-only values (e.g. the
        # "Antal" column's "1" / "5") are not re-parsed as real numbers,
        # then strip the number-format style again afterwards so the
        # cell is left as plain, unstyled text - matching how the sheet
        # already stores every other text column.
        $targetCell.NumberFormat = "@"
        $targetCell.Value2 = $newValue
        $targetCell.ClearFormats()
    } else {
        $targetCell.Value2 = $newValue
    }
}

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell13 = $ws.Cells.Item(13, $c)
    $cell14 = $ws.Cells.Item(14, $c)

    $v13 = $cell13.Value2
    $v14 = $cell14.Value2

    # Only touch cells whose content actually differs between the two
    # rows - this avoids re-typing identical values (e.g. matching dates)
    # back into a cell, which would otherwise make Excel re-infer the
    # cell type/format (turning a plain date-text string into a real
    # date serial number) even though nothing really changed.
    if ($v13 -eq $v14) {
        continue
    }

    $isNumeric = $numericCols -contains $c
    $isBoolean = $booleanCols -contains $c
    $isText = -not $isNumeric -and -not $isBoolean

    Set-SwappedValue $cell13 $v14 $isText
    Set-SwappedValue $cell14 $v13 $isText
}
